# Auto-generated data-driven update of cryptos list (sheet1) cell values.
# Mirrors the authoritative unified diff: updates Price (D) and Volume(1h) (E)
# columns for the affected rows, plus the swapped ARBITRUM / MXToken rows (39/40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-seven character used in the BabyDogeCoin price (row 49), e.g. 0.0₇0974
$sub7 = [char]0x2087

# Each entry: row number plus the new values for any of the B/C/D/E columns that changed.
# ForceText marks Price values that look numeric (e.g. "212.33") so Excel's auto-
# number detection doesn't convert them away from text, matching the source cell type.
$updates = @(
    @{Row=2; D='26.190.00'; E='  -0.30%  '},
    @{Row=3; D='1.599.70'; E='  -0.63%  '},
    @{Row=4; E='  +0.07%  '},
    @{Row=5; D='212.33'; E='  -0.32%  '; ForceText=$true},
    @{Row=6; E='  +0.07%  '},
    @{Row=7; E='  +0.25%  '},
    @{Row=8; E='  -0.35%  '},
    @{Row=9; E='  -1.02%  '},
    @{Row=10; E='  -0.29%  '},
    @{Row=11; E='  -0.56%  '},
    @{Row=12; D='1.823.50'; E='  -0.44%  '},
    @{Row=13; D='1.594.35'; E='  -0.83%  '},
    @{Row=14; D='4.02'; E='  +0.05%  '; ForceText=$true},
    @{Row=15; D='0.513'; E='  +0.36%  '; ForceText=$true},
    @{Row=16; D='26.173.41'; E='  -0.28%  '},
    @{Row=17; D='61.28'; E='  +0.78%  '; ForceText=$true},
    @{Row=18; E='  -0.17%  '},
    @{Row=19; E='  +0.06%  '},
    @{Row=20; D='201.81'; E='  +1.03%  '; ForceText=$true},
    @{Row=21; D='4.28'; E='  +0.77%  '; ForceText=$true},
    @{Row=22; D='9.24'; E='  -1.69%  '; ForceText=$true},
    @{Row=23; D='5.96'; E='  -1.09%  '; ForceText=$true},
    @{Row=24; D='1.92'; E='  +9.12%  '; ForceText=$true},
    @{Row=25; D='144.04'; E='  +0.92%  '; ForceText=$true},
    @{Row=26; E='  +0.13%  '},
    @{Row=27; E='  -7.78%  '},
    @{Row=28; D='15.15'; E='  -0.33%  '; ForceText=$true},
    @{Row=29; D='6.55'; E='  +0.67%  '; ForceText=$true},
    @{Row=30; E='  +3.02%  '},
    @{Row=31; E='  -0.45%  '},
    @{Row=32; D='3.18'; E='  +1.60%  '; ForceText=$true},
    @{Row=33; E='  -3.17%  '},
    @{Row=34; E='  +3.19%  '},
    @{Row=35; E='  -0.80%  '},
    @{Row=36; D='1.151.90'; E='  +3.70%  '},
    @{Row=37; E='  +7.46%  '},
    @{Row=38; E='  +0.11%  '},
    @{Row=39; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='0.792'; E='  +0.74%  '; ForceText=$true},
    @{Row=40; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.33'; E='  -1.18%  '; ForceText=$true},
    @{Row=41; D='0.497'; E='  -1.23%  '; ForceText=$true},
    @{Row=42; D='0.780'; E='  +0.15%  '; ForceText=$true},
    @{Row=43; E='  +2.02%  '},
    @{Row=44; D='1.737.22'; E='  -0.43%  '},
    @{Row=45; D='91.65'; ForceText=$true},
    @{Row=46; E='  -3.18%  '},
    @{Row=47; D='53.98'; E='  +0.24%  '; ForceText=$true},
    @{Row=48; E='  -0.84%  '},
    @{Row=49; D=("0.0{0}0974" -f $sub7); E='  -8.57%  '},
    @{Row=50; E='  -0.70%  '},
    @{Row=51; E='  -0.08%  '}
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($item.ContainsKey("B")) { $ws.Range("B" + $r).Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Range("C" + $r).Value = $item.C }
    if ($item.ContainsKey("D")) {
        if ($item.ContainsKey("ForceText")) {
            # Force text storage so numeric-looking prices are not auto-converted
            # to numbers by Excel, matching the original inline-string cells.
            $ws.Range("D" + $r).NumberFormat = "@"
        }
        $ws.Range("D" + $r).Value = $item.D
    }
    if ($item.ContainsKey("E")) { $ws.Range("E" + $r).Value = $item.E }
}

Write-Output "Applied $($updates.Count) row updates."
